{"js": "// Update the date line and the 25 division problems in the practice table.\n// Each old value is unique in the document, but some NEW values collide\n// with OLD values elsewhere in the table (e.g. \"81\u00f77=\" is the old text of\n// one cell and the new text of another). To stay unambiguous we scope each\n// search to the specific paragraph / table cell it belongs to, rather than\n// searching the whole document body, and we rewrite the found range in\n// place (InsertLocation.Replace on the matched range) so the existing run\n// formatting (font/size) and paragraph formatting (alignment) are kept.\n\n// 1) Title / date line: first paragraph of the body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst titleResults = titlePara.search(\"2024-04-28 Sunday\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the date line, found \" + titleResults.items.length);\n}\ntitleResults.items[0].insertText(\"2024-04-29 Monday\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) The division problems table (5 columns; only rows 0, 4, 8, 12, 16 hold text).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"91\u00f75=\", newText: \"92\u00f72=\" },\n  { row: 0, col: 1, oldText: \"86\u00f76=\", newText: \"81\u00f77=\" },\n  { row: 0, col: 2, oldText: \"14\u00f72=\", newText: \"15\u00f78=\" },\n  { row: 0, col: 3, oldText: \"37\u00f75=\", newText: \"15\u00f74=\" },\n  { row: 0, col: 4, oldText: \"19\u00f75=\", newText: \"93\u00f74=\" },\n\n  { row: 4, col: 0, oldText: \"94\u00f74=\", newText: \"75\u00f78=\" },\n  { row: 4, col: 1, oldText: \"17\u00f77=\", newText: \"35\u00f72=\" },\n  { row: 4, col: 2, oldText: \"27\u00f72=\", newText: \"70\u00f77=\" },\n  { row: 4, col: 3, oldText: \"16\u00f78=\", newText: \"22\u00f73=\" },\n  { row: 4, col: 4, oldText: \"10\u00f72=\", newText: \"13\u00f78=\" },\n\n  { row: 8, col: 0, oldText: \"67\u00f76=\", newText: \"44\u00f74=\" },\n  { row: 8, col: 1, oldText: \"60\u00f75=\", newText: \"78\u00f79=\" },\n  { row: 8, col: 2, oldText: \"82\u00f77=\", newText: \"83\u00f79=\" },\n  { row: 8, col: 3, oldText: \"20\u00f75=\", newText: \"62\u00f76=\" },\n  { row: 8, col: 4, oldText: \"81\u00f77=\", newText: \"43\u00f76=\" },\n\n  { row: 12, col: 0, oldText: \"98\u00f77=\", newText: \"90\u00f79=\" },\n  { row: 12, col: 1, oldText: \"77\u00f74=\", newText: \"99\u00f78=\" },\n  { row: 12, col: 2, oldText: \"38\u00f73=\", newText: \"64\u00f76=\" },\n  { row: 12, col: 3, oldText: \"71\u00f72=\", newText: \"23\u00f76=\" },\n  { row: 12, col: 4, oldText: \"16\u00f75=\", newText: \"46\u00f72=\" },\n\n  { row: 16, col: 0, oldText: \"46\u00f79=\", newText: \"74\u00f72=\" },\n  { row: 16, col: 1, oldText: \"76\u00f72=\", newText: \"40\u00f77=\" },\n  { row: 16, col: 2, oldText: \"63\u00f75=\", newText: \"82\u00f73=\" },\n  { row: 16, col: 3, oldText: \"49\u00f78=\", newText: \"44\u00f75=\" },\n  { row: 16, col: 4, oldText: \"54\u00f74=\", newText: \"42\u00f75=\" },\n];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const found = cell.body.search(r.oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly one match for '\" + r.oldText + \"' in cell (\" + r.row + \",\" + r.col + \"), found \" + found.items.length\n    );\n  }\n  found.items[0].insertText(r.newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice table.\n#\n# NOTE: this host's Find.Execute(..., Replace:=wdReplaceAll) ignores the\n# bounding Range it was invoked on and replaces matches anywhere in the\n# document \u2014 so scoping Find to a single cell's Range is not safe here.\n# That matters because some NEW values collide with OLD values elsewhere\n# in the table (e.g. \"81\u00f77=\" is the old text of one cell and the new text\n# of another): a document-wide Find/Replace run in sequence would corrupt\n# a cell that was already updated. Instead we target each paragraph/cell\n# Range directly and assign its .Text property, which only touches that\n# Range and leaves the existing run formatting (font/size) and paragraph\n# formatting (alignment) of that run in place.\n#\n# (A table cell's Range.Text carries trailing cell-mark control characters,\n# so we verify the old value with StartsWith rather than an exact -eq.)\n\n$d = $word.ActiveDocument\n\n# 1) Title / date line: first paragraph of the body.\n$titlePara = $d.Paragraphs.Item(1)\nif (-not $titlePara.Range.Text.StartsWith(\"2024-04-28 Sunday\")) {\n    throw \"Title paragraph did not contain the expected date text.\"\n}\n$titlePara.Range.Text = \"2024-04-29 Monday\"\n\n# 2) The division problems table (5 columns; only rows 1, 5, 9, 13, 17 hold text;\n#    Word COM Cell()/Paragraphs indices are 1-based).\n$tbl = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"91\u00f75=\"; New = \"92\u00f72=\" },\n    @{ Row = 1;  Col = 2; Old = \"86\u00f76=\"; New = \"81\u00f77=\" },\n    @{ Row = 1;  Col = 3; Old = \"14\u00f72=\"; New = \"15\u00f78=\" },\n    @{ Row = 1;  Col = 4; Old = \"37\u00f75=\"; New = \"15\u00f74=\" },\n    @{ Row = 1;  Col = 5; Old = \"19\u00f75=\"; New = \"93\u00f74=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"94\u00f74=\"; New = \"75\u00f78=\" },\n    @{ Row = 5;  Col = 2; Old = \"17\u00f77=\"; New = \"35\u00f72=\" },\n    @{ Row = 5;  Col = 3; Old = \"27\u00f72=\"; New = \"70\u00f77=\" },\n    @{ Row = 5;  Col = 4; Old = \"16\u00f78=\"; New = \"22\u00f73=\" },\n    @{ Row = 5;  Col = 5; Old = \"10\u00f72=\"; New = \"13\u00f78=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"67\u00f76=\"; New = \"44\u00f74=\" },\n    @{ Row = 9;  Col = 2; Old = \"60\u00f75=\"; New = \"78\u00f79=\" },\n    @{ Row = 9;  Col = 3; Old = \"82\u00f77=\"; New = \"83\u00f79=\" },\n    @{ Row = 9;  Col = 4; Old = \"20\u00f75=\"; New = \"62\u00f76=\" },\n    @{ Row = 9;  Col = 5; Old = \"81\u00f77=\"; New = \"43\u00f76=\" },\n\n    @{ Row = 13; Col = 1; Old = \"98\u00f77=\"; New = \"90\u00f79=\" },\n    @{ Row = 13; Col = 2; Old = \"77\u00f74=\"; New = \"99\u00f78=\" },\n    @{ Row = 13; Col = 3; Old = \"38\u00f73=\"; New = \"64\u00f76=\" },\n    @{ Row = 13; Col = 4; Old = \"71\u00f72=\"; New = \"23\u00f76=\" },\n    @{ Row = 13; Col = 5; Old = \"16\u00f75=\"; New = \"46\u00f72=\" },\n\n    @{ Row = 17; Col = 1; Old = \"46\u00f79=\"; New = \"74\u00f72=\" },\n    @{ Row = 17; Col = 2; Old = \"76\u00f72=\"; New = \"40\u00f77=\" },\n    @{ Row = 17; Col = 3; Old = \"63\u00f75=\"; New = \"82\u00f73=\" },\n    @{ Row = 17; Col = 4; Old = \"49\u00f78=\"; New = \"44\u00f75=\" },\n    @{ Row = 17; Col = 5; Old = \"54\u00f74=\"; New = \"42\u00f75=\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $tbl.Cell($r.Row, $r.Col)\n    if (-not $cell.Range.Text.StartsWith($r.Old)) {\n        throw (\"Cell (\" + $r.Row + \",\" + $r.Col + \") did not contain the expected text '\" + $r.Old + \"'.\")\n    }\n    $cell.Range.Text = $r.New\n}\n"}
